$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row (A1:D1) from Spanish labels to snake_case field names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2) Title-case the connector words (de/del/la/el/los/las/y) inside municipality /
#    state names, e.g. "Pabellón de Arteaga" -> "Pabellón De Arteaga"
$titleCaseFixes = @(
  @("B8","Pabellón De Arteaga"),
  @("B9","Rincón De Romos"),
  @("B10","San Francisco De Los Romo"),
  @("B11","San José De Gracia"),
  @("B16","Playas De Rosarito"),
  @("B41","Amatenango De La Frontera"),
  @("B44","Bejucal De Ocampo"),
  @("B46","Benemérito De Las Américas"),
  @("B55","Chiapa De Corzo"),
  @("B60","Comitán De Domínguez"),
  @("B83","Marqués De Comillas"),
  @("B86","Montecristo De Guerrero"),
  @("B91","Ocozocoautla De Espinosa"),
  @("B99","Salto De Agua"),
  @("B100","San Cristóbal De Las Casas"),
  @("B144","Coyame Del Sotol"),
  @("B155","Guadalupe Y Calvo"),
  @("B158","Hidalgo Del Parral"),
  @("B184","San Francisco De Borja"),
  @("B185","San Francisco De Conchos"),
  @("B186","San Francisco Del Oro"),
  @("B194","Valle De Zaragoza"),
  @("B227","San Juan De Sabinas"),
  @("B244","Villa De Álvarez"),
  @("A246","Ciudad De México"),
  @("B250","Cuajimalpa De Morelos"),
  @("B265","Coneto De Comonfort"),
  @("B279","Nombre De Dios"),
  @("B283","Pánuco De Coronado"),
  @("B290","San Juan De Guadalupe"),
  @("B291","San Juan Del Río"),
  @("B292","San Luis Del Cordero"),
  @("B293","San Pedro Del Gallo"),
  @("A303","Estado De México"),
  @("B303","Acambay De Ruíz Castañeda"),
  @("B306","Almoloya De Alquisiras"),
  @("B307","Almoloya De Juárez"),
  @("B308","Almoloya Del Río"),
  @("B314","Atizapán De Zaragoza"),
  @("B322","Chapa De Mota"),
  @("B326","Coacalco De Berriozábal"),
  @("B333","Ecatepec De Morelos"),
  @("B341","Ixtapan De La Sal"),
  @("B342","Ixtapan Del Oro"),
  @("B358","Naucalpan De Juárez"),
  @("B370","San Antonio La Isla"),
  @("B371","San Felipe Del Progreso"),
  @("B372","San Martín De Las Pirámides"),
  @("B374","San Simón De Guerrero"),
  @("B376","Soyaniquilpan De Juárez"),
  @("B386","Tenango Del Aire"),
  @("B387","Tenango Del Valle"),
  @("B400","Tlalnepantla De Baz"),
  @("B406","Valle De Bravo"),
  @("B407","Valle De Chalco Solidaridad"),
  @("B408","Villa De Allende"),
  @("B409","Villa Del Carbón"),
  @("B422","San Miguel De Allende"),
  @("B423","Apaseo El Alto"),
  @("B424","Apaseo El Grande"),
  @("B432","Dolores Hidalgo Cuna De La Independencia Nacional"),
  @("B436","Jaral Del Progreso"),
  @("B444","Purísima Del Rincón"),
  @("B448","San Diego De La Unión"),
  @("B450","San Francisco Del Rincón"),
  @("B452","San Luis De La Paz"),
  @("B454","Santa Cruz De Juventino Rosas"),
  @("B456","Silao De La Victoria"),
  @("B461","Valle De Santiago"),
  @("B467","Acapulco De Juárez"),
  @("B470","Ajuchitlán Del Progreso"),
  @("B471","Alcozauca De Guerrero"),
  @("B475","Atenango Del Río"),
  @("B476","Atlamajalcingo Del Monte"),
  @("B478","Atoyac De Álvarez"),
  @("B479","Ayutla De Los Libres"),
  @("B482","Buenavista De Cuéllar"),
  @("B483","Chilapa De Álvarez"),
  @("B484","Chilpancingo De Los Bravo"),
  @("B485","Coahuayutla De José María Izazaga"),
  @("B490","Coyuca De Benítez"),
  @("B491","Coyuca De Catalán"),
  @("B495","Cuetzala Del Progreso"),
  @("B496","Cutzamala De Pinzón"),
  @("B502","Huitzuco De Los Figueroa"),
  @("B503","Iguala De La Independencia"),
  @("B505","Ixcateopan De Cuauhtémoc"),
  @("B506","Zihuatanejo De Azueta"),
  @("B508","La Unión De Isidoro Montes De Oca"),
  @("B511","Mártir De Cuilapan"),
  @("B524","Taxco De Alarcón"),
  @("B526","Técpan De Galeana"),
  @("B528","Tepecoacuilco De Trujano"),
  @("B530","Tixtla De Guerrero"),
  @("B534","Tlalixtaquilla De Maldonado"),
  @("B535","Tlapa De Comonfort"),
  @("B547","Agua Blanca De Iturbide"),
  @("B554","Atotonilco De Tula"),
  @("B555","Atotonilco El Grande"),
  @("B561","Cuautepec De Hinojosa"),
  @("B567","Huasca De Ocampo"),
  @("B571","Huejutla De Reyes"),
  @("B574","Jacala De Ledezma"),
  @("B581","Mineral De La Reforma"),
  @("B582","Mineral Del Chico"),
  @("B583","Mineral Del Monte"),
  @("B584","Mixquiahuala De Juárez"),
  @("B585","Molango De Escamilla"),
  @("B587","Nopala De Villagrán"),
  @("B588","Omitlán De Juárez"),
  @("B589","Pachuca De Soto"),
  @("B592","Progreso De Obregón"),
  @("B598","Santiago De Anaya"),
  @("B599","Santiago Tulantepec De Lugo Guerrero"),
  @("B603","Tenango De Doria"),
  @("B605","Tepehuacán De Guerrero"),
  @("B606","Tepeji Del Río De Ocampo"),
  @("B609","Tezontepec De Aldama"),
  @("B618","Tula De Allende"),
  @("B619","Tulancingo De Bravo"),
  @("B620","Villa De Tezontepec"),
  @("B624","Zacualtipán De Ángeles"),
  @("B625","Zapotlán De Juárez"),
  @("B630","Acatlán De Juárez"),
  @("B631","Ahualulco De Mercado"),
  @("B636","Atemajac De Brizuela"),
  @("B639","Atotonilco El Alto"),
  @("B641","Autlán De Navarro"),
  @("B653","Concepción De Buenos Aires"),
  @("B654","Cuautitlán De García Barragán"),
  @("B662","Encarnación De Díaz"),
  @("B669","Huejuquilla El Alto"),
  @("B670","Ixtlahuacán De Los Membrillos"),
  @("B671","Ixtlahuacán Del Río"),
  @("B675","Jilotlán De Los Dolores"),
  @("B681","La Manzanilla De La Paz"),
  @("B682","Lagos De Moreno"),
  @("B690","Ojuelos De Jalisco"),
  @("B695","San Cristóbal De La Barranca"),
  @("B696","San Diego De Alejandría"),
  @("B698","San Juan De Los Lagos"),
  @("B699","San Juanito De Escobedo"),
  @("B702","San Martín De Bolaños"),
  @("B704","San Miguel El Alto"),
  @("B705","San Sebastián Del Oeste"),
  @("B706","Santa María De Los Ángeles"),
  @("B707","Santa María Del Oro"),
  @("B710","Talpa De Allende"),
  @("B711","Tamazula De Gordiano"),
  @("B714","Techaluta De Montenegro"),
  @("B718","Teocuitatlán De Corona"),
  @("B719","Tepatitlán De Morelos"),
  @("B722","Tizapán El Alto"),
  @("B723","Tlajomulco De Zúñiga"),
  @("B735","Unión De San Antonio"),
  @("B736","Unión De Tula"),
  @("B737","Valle De Guadalupe"),
  @("B738","Valle De Juárez"),
  @("B743","Yahualica De González Gallo"),
  @("B744","Zacoalco De Torres"),
  @("B747","Zapotitlán De Vadillo"),
  @("B748","Zapotlán Del Rey"),
  @("B749","Zapotlán El Grande"),
  @("B775","Coalcomán De Vázquez Pallares"),
  @("B777","Cojumatlán De Régules"),
  @("B844","Tiquicheo De Nicolás Romero"),
  @("B870","Coatlán Del Río"),
  @("B878","Jonacatepec De Leandro Valle"),
  @("B882","Puente De Ixtla"),
  @("B888","Tetela Del Volcán"),
  @("B889","Tlaltizapán De Zapata"),
  @("B897","Zacualpan De Amilpas"),
  @("B901","Amatlán De Cañas"),
  @("B902","Bahía De Banderas"),
  @("B905","Ixtlán Del Río"),
  @("B912","Santa María Del Oro"),
  @("B930","Ciénega De Flores"),
  @("B948","Lampazos De Naranjo"),
  @("B955","Mier Y Noriega"),
  @("B964","San Nicolás De Los Garza"),
  @("B972","Acatlán De Pérez Figueroa"),
  @("B978","Ayoquezco De Aldama"),
  @("B983","Chalcatongo De Hidalgo"),
  @("B984","Ciénega De Zimatlán"),
  @("B987","Coicoyán De Las Flores"),
  @("B989","Constancia Del Rosario"),
  @("B991","Cuilápam De Guerrero"),
  @("B992","El Barrio De La Soledad"),
  @("B993","Fresnillo De Trujano"),
  @("B994","Guadalupe De Ramírez"),
  @("B996","Guevea De Humboldt"),
  @("B997","Heroica Ciudad De Ejutla De Crespo"),
  @("B998","Heroica Ciudad De Huajuapan De León"),
  @("B999","Heroica Ciudad De Tlaxiaco"),
  @("B1001","Huautla De Jiménez"),
  @("B1002","Ixtlán De Juárez"),
  @("B1003","Heroica Ciudad De Juchitán De Zaragoza"),
  @("B1012","Mariscala De Juárez"),
  @("B1015","Miahuatlán De Porfirio Díaz"),
  @("B1016","Mixistlán De La Reforma"),
  @("B1019","Nejapa De Madero"),
  @("B1020","Oaxaca De Juárez"),
  @("B1021","Ocotlán De Morelos"),
  @("B1022","Pinotepa De Don Luis"),
  @("B1024","Putla Villa De Guerrero"),
  @("B1037","San Antonino El Alto"),
  @("B1051","San Dionisio Del Mar"),
  @("B1054","San Felipe Jalapa De Díaz"),
  @("B1057","San Francisco Del Mar"),
  @("B1076","San José Del Progreso"),
  @("B1091","San Juan De Los Cués"),
  @("B1092","San Juan Del Estado"),
  @("B1093","San Juan Del Río"),
  @("B1120","San Martín De Los Cansecos"),
  @("B1133","San Miguel Del Puerto"),
  @("B1135","San Miguel El Grande"),
  @("B1147","San Pablo Villa De Mitla"),
  @("B1152","San Pedro El Alto"),
  @("B1168","San Pedro Y San Pablo Ayutla"),
  @("B1169","San Pedro Y San Pablo Teposcolula"),
  @("B1170","San Pedro Y San Pablo Tequixtepec"),
  @("B1178","Santa Ana Del Valle"),
  @("B1187","Santa Cruz Tacache De Mina"),
  @("B1191","Santa Inés De Zaragoza"),
  @("B1192","Santa Inés Del Monte"),
  @("B1194","Santa Lucía Del Camino"),
  @("B1205","Santa María Jalapa Del Marqués"),
  @("B1227","Santiago Del Río"),
  @("B1257","Santo Domingo De Morelos"),
  @("B1274","Tamazulápam Del Espíritu Santo"),
  @("B1275","Tataltepec De Valdés"),
  @("B1276","Teotitlán De Flores Magón"),
  @("B1277","Teotitlán Del Valle"),
  @("B1278","Tepelmeme Villa De Morelos"),
  @("B1279","Tezoatlán De Segura Y Luna"),
  @("B1280","Tlacolula De Matamoros"),
  @("B1281","Tlalixtac De Cabrera"),
  @("B1282","Totontepec Villa De Morelos"),
  @("B1284","Villa De Chilapa De Díaz"),
  @("B1285","Villa De Etla"),
  @("B1286","Villa De Tututepec De Melchor Ocampo"),
  @("B1287","Villa De Zaachila"),
  @("B1290","Villa Sola De Vega"),
  @("B1291","Villa Talea De Castro"),
  @("B1292","Villa Tejúpam De La Unión"),
  @("B1294","Zapotitlán Del Río"),
  @("B1296","Zimatlán De Álvarez"),
  @("B1320","Ayotoxco De Guerrero"),
  @("B1324","Chalchicomula De Sesma"),
  @("B1334","Chila De La Sal"),
  @("B1345","Cuayuca De Andrade"),
  @("B1346","Cuetzalan Del Progreso"),
  @("B1362","Huehuetlán El Chico"),
  @("B1363","Huehuetlán El Grande"),
  @("B1368","Huitzilan De Serdán"),
  @("B1370","Ixcamilpa De Guerrero"),
  @("B1374","Izúcar De Matamoros"),
  @("B1384","Los Reyes De Juárez"),
  @("B1392","Palmar De Bravo"),
  @("B1401","San Diego La Mesa Tochimiltzingo"),
  @("B1414","San Nicolás De Los Ranchos"),
  @("B1418","San Salvador El Seco"),
  @("B1419","San Salvador El Verde"),
  @("B1426","Tecali De Herrera"),
  @("B1434","Tepanco De López"),
  @("B1435","Tepango De Rodríguez"),
  @("B1436","Tepatlaxco De Hidalgo"),
  @("B1441","Tepexi De Rodríguez"),
  @("B1443","Tepeyahualco De Cuauhtémoc"),
  @("B1444","Tetela De Ocampo"),
  @("B1445","Teteles De Avila Castillo"),
  @("B1450","Tlacotepec De Benito Juárez"),
  @("B1462","Totoltepec De Guerrero"),
  @("B1464","Tuzamapan De Galeana"),
  @("B1468","Xayacatlán De Bravo"),
  @("B1474","Xochitlán De Vicente Suárez"),
  @("B1482","Zapotitlán De Méndez"),
  @("B1491","Amealco De Bonfil"),
  @("B1493","Cadereyta De Montes"),
  @("B1499","Jalpan De Serra"),
  @("B1500","Landa De Matamoros"),
  @("B1503","Pinal De Amoles"),
  @("B1506","San Juan Del Río"),
  @("B1521","Armadillo De Los Infante"),
  @("B1522","Axtla De Terrazas"),
  @("B1527","Cerro De San Pedro"),
  @("B1529","Ciudad Del Maíz"),
  @("B1540","Mexquitic De Carmona"),
  @("B1546","San Ciro De Acosta"),
  @("B1552","Santa María Del Río"),
  @("B1554","Soledad De Graciano Sánchez"),
  @("B1562","Tanquián De Escobedo"),
  @("B1566","Villa De Arista"),
  @("B1567","Villa De Arriaga"),
  @("B1568","Villa De Guadalupe"),
  @("B1569","Villa De La Paz"),
  @("B1570","Villa De Ramos"),
  @("B1571","Villa De Reyes"),
  @("B1618","Nacozari De García"),
  @("B1625","San Miguel De Horcasitas"),
  @("B1637","Jalpa De Méndez"),
  @("B1682","Soto La Marina"),
  @("B1690","Acuamanala De Miguel Hidalgo"),
  @("B1692","Amaxac De Guerrero"),
  @("B1693","Apetatitlán De Antonio Carvajal"),
  @("B1698","Contla De Juan Cuamatzi"),
  @("B1704","Ixtacuixtla De Mariano Matamoros"),
  @("B1707","Mazatecochco De José María Morelos"),
  @("B1708","Nanacamilpa De Mariano Arista"),
  @("B1711","Papalotla De Xicohténcatl"),
  @("B1714","San Pablo Del Monte"),
  @("B1715","Sanctórum De Lázaro Cárdenas"),
  @("B1719","Tepetitla De Lardizábal"),
  @("B1722","Tetla De La Solidaridad"),
  @("B1734","Ziltlaltépec De Trinidad Sánchez Santos"),
  @("B1744","Alto Lucero De Gutiérrez Barrios"),
  @("B1748","Amatlán De Los Reyes"),
  @("B1760","Boca Del Río"),
  @("B1762","Camarón De Tejeda"),
  @("B1766","Castillo De Teayo"),
  @("B1768","Cazones De Herrera"),
  @("B1776","Chinampa De Gorostiza"),
  @("B1789","Cosamaloapan De Carpio"),
  @("B1790","Cosautlán De Carvajal"),
  @("B1807","Hueyapan De Ocampo"),
  @("B1808","Ignacio De La Llave"),
  @("B1812","Ixhuacán De Los Reyes"),
  @("B1813","Ixhuatlán De Madero"),
  @("B1814","Ixhuatlán Del Café"),
  @("B1815","Ixhuatlán Del Sureste"),
  @("B1826","Juchique De Ferrer"),
  @("B1829","Landero Y Coss"),
  @("B1832","Las Vigas De Ramírez"),
  @("B1833","Lerdo De Tejada"),
  @("B1837","Martínez De La Torre"),
  @("B1840","Medellín De Bravo"),
  @("B1845","Nanchital De Lázaro Cárdenas Del Río"),
  @("B1856","Ozuluama De Mascareñas"),
  @("B1860","Paso De Ovejas"),
  @("B1861","Paso Del Macho"),
  @("B1865","Poza Rica De Hidalgo"),
  @("B1876","Sayula De Alemán"),
  @("B1880","Soledad De Doblado"),
  @("B1888","Tatahuicapan De Juárez"),
  @("B1910","Tlacotepec De Mejía"),
  @("B1924","Vega De Alatorre"),
  @("B1935","Zontecomatlán De López Y Fuentes"),
  @("B1936","Zozocolco De Hidalgo"),
  @("B1977","Cañitas De Felipe Pescador"),
  @("B1979","Concepción Del Oro"),
  @("B1981","El Plateado De Joaquín Amaro"),
  @("B1992","Jiménez Del Teul"),
  @("B1999","Mezquital Del Oro"),
  @("B2004","Moyahua De Estrada"),
  @("B2005","Nochistlán De Mejía"),
  @("B2006","Noria De Ángeles"),
  @("B2016","Teúl De González Ortega"),
  @("B2017","Tlaltenango De Sánchez Román"),
  @("B2019","Trinidad García De La Cadena"),
  @("B2022","Villa De Cos")
)
foreach ($fix in $titleCaseFixes) {
  $ws.Range($fix[0]).Value = $fix[1]
}

# 3) A handful of percentage cells had their floating point literal re-rounded
#    by one ULP when the source file was last resaved by Excel; reproduce the
#    exact target values.
$floatFixes = @(
  @("D57",0.0000920132499079867621),
  @("D97",0.0000920132499079867621),
  @("D140",0.0000920132499079867621),
  @("D151",0.0000920132499079867621),
  @("D155",0.0000920132499079867621),
  @("D158",0.0009488866396761132076),
  @("D199",0.0000920132499079867621),
  @("D213",0.0094486105999263877292),
  @("D304",0.0000920132499079867621),
  @("D348",0.0000920132499079867621),
  @("D416",0.0000920132499079867621),
  @("D449",0.0090748067721751913411),
  @("D452",0.0098971751932278237418),
  @("D526",0.0009431358115568640307),
  @("D566",0.0000920132499079867621),
  @("D585",0.0000920132499079867621),
  @("D594",0.0000920132499079867621),
  @("D631",0.0000920132499079867621),
  @("D671",0.0000920132499079867621),
  @("D680",0.0000920132499079867621),
  @("D717",0.0009603882959146115613),
  @("D738",0.0000920132499079867621),
  @("D886",0.0000920132499079867621),
  @("D936",0.0000920132499079867621),
  @("D940",0.0000920132499079867621),
  @("D961",0.0000920132499079867621),
  @("D1298",0.0000920132499079867621),
  @("D1352",0.0000920132499079867621),
  @("D1374",0.0009603882959146115613),
  @("D1520",0.0009431358115568640307),
  @("D1536",0.0000920132499079867621),
  @("D1586",0.0000920132499079867621),
  @("D1673",0.0000920132499079867621),
  @("D1730",0.0000920132499079867621),
  @("D1793",0.0000920132499079867621),
  @("D1867",0.0000920132499079867621),
  @("D1870",0.0000920132499079867621),
  @("D1936",0.0000920132499079867621),
  @("D1971",0.0009431358115568640307),
  @("D2025",0.0009201324990798676485)
)
foreach ($fix in $floatFixes) {
  $ws.Range($fix[0]).Value = $fix[1]
}

# 4) Drop the trailing footnote rows (sample size / source / author / date)
#    that used to sit below the data table.
$ws.Range("A2031:A2035").EntireRow.Delete()
